$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 72 (existing rows 72:136 shift down to 73:137).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "Vega Modelo de Temuco"
$ws.Range("C72").Value = "La Araucanía"
$ws.Range("D72").Value = 45280
$ws.Range("E72").Value = 9
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100108
$ws.Range("H72").Value = "Tropicales y subtropicales"
$ws.Range("I72").Value = 100108007
$ws.Range("J72").Value = "Coco"
$ws.Range("K72").Value = "Sin especificar"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 60
$ws.Range("N72").Value = 32000
$ws.Range("O72").Value = 32000
$ws.Range("P72").Value = 32000
$ws.Range("Q72").Value = "$/malla 20 unidades"
$ws.Range("R72").Value = "Perú"
$ws.Range("S72").Value = 1600
$ws.Range("T72").Value = 20
